# edit.ps1 — apply the two textual corrections described by the diff:
#
#  1. "Et leur reste sur la premiere annee quatre vintg livres quinze"
#     gains a trailing space and is followed by a new run containing
#     "[sols]" (minimal run formatting: only <w:rtl w:val="0"/>).
#
#  2. The run containing "devant." loses its trailing period, becoming
#     "devant" (same run formatting retained).

$d = $word.ActiveDocument

# --- Change 1a: add the trailing space to the "quinze" run -----------------
$d.Content.Find.Execute(
    "quatre vintg livres quinze", $true, $false, $false, $false, $false,
    $true, 1, $false, "quatre vintg livres quinze ", 2) | Out-Null

# --- Change 1b: insert a new "[sols]" run right after it, before "<lb/>" ---
# We need the new run's rPr to be just <w:rtl w:val="0"/> (no rFonts/color/
# size), matching the paragraph's existing convention for "plain" runs
# (e.g. the single-letter runs used elsewhere in this document). The
# cleanest way to reproduce that exact formatting via the object model is
# to clone it from Range.FormattedText of an existing run that already has
# that minimal formatting, then overwrite the copied text.
$full = $d.Content.Text

$srcMarker = "Transport aulx heritiers "
$srcStart = $full.IndexOf($srcMarker) + $srcMarker.Length
$srcRange = $d.Range($srcStart, $srcStart + 1)

$marker = "quatre vintg livres quinze "
$insertPos = $full.IndexOf($marker) + $marker.Length
$insertionPoint = $d.Range($insertPos, $insertPos)
$insertionPoint.FormattedText = $srcRange.FormattedText

$newRun = $d.Range($insertPos, $insertPos + 1)
$newRun.Text = "[sols]"

# --- Change 2: "devant." -> "devant" ---------------------------------------
$d.Content.Find.Execute(
    "devant.", $true, $false, $false, $false, $false,
    $true, 1, $false, "devant", 2) | Out-Null

Write-Output "done"
